$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.35249999999999
$ws.Range("C21").Value = -13.38220000000001
$ws.Range("C23").Value = -12.1345
$ws.Range("D24").Value = -7.641699999999998
$ws.Range("C25").Value = -11.65040000000001
$ws.Range("D28").Value = -8.327499999999993
$ws.Range("D36").Value = -7.7137
$ws.Range("D45").Value = -7.231099999999998
$ws.Range("D48").Value = -7.748699999999994
$ws.Range("D49").Value = -8.078300000000002
$ws.Range("D52").Value = -7.927700000000003
$ws.Range("C53").Value = -11.8085
$ws.Range("D53").Value = -8.086599999999994
$ws.Range("D54").Value = -7.924999999999998
$ws.Range("C57").Value = -14.21489999999999
$ws.Range("C59").Value = -13.05930000000001
$ws.Range("C69").Value = -10.7491
$ws.Range("D70").Value = -6.9779
$ws.Range("C79").Value = -11.93620000000001
$ws.Range("C83").Value = -13.49889999999999
$ws.Range("D86").Value = -8.303599999999998
$ws.Range("D87").Value = -8.273999999999992
$ws.Range("C93").Value = -10.32509999999999
$ws.Range("D101").Value = -7.978500000000002
